$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Duplicate the "2017" sheet to create "2018" BEFORE editing "2017" so the
#    new sheet inherits the original (pre-edit) formulas/values - matches the
#    new xl/worksheets/sheet4.xml which still has the 'Initial Buys'! F-column
#    formulas and the original D17/S17 values.
# ---------------------------------------------------------------------------
$ws2017 = $wb.Worksheets.Item("2017")
$ws2017.Copy($null, $ws2017) | Out-Null

$lastIndex = $wb.Worksheets.Count
$ws2018 = $wb.Worksheets.Item($lastIndex)
$ws2018.Name = "2018"

# ---------------------------------------------------------------------------
# 2. On the new "2018" sheet: point the running-dividend-total column (G) at
#    last year's totals ('2017'!U..) instead of '2016'!S.., and zero out the
#    per-month dividend columns (H:S) since no 2018 dividends have posted yet
#    -- except January (H) on row 17, which already received one payment.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 18; $r++) {
    $ws2018.Range("G$r").Formula = "='2017'!U$r"
    foreach ($col in @("H","I","J","K","L","M","N","O","P","Q","R","S")) {
        $ws2018.Range("$col$r").Value = 0
    }
}
$ws2018.Range("H17").Value = 7.92

# ---------------------------------------------------------------------------
# 3. Back on "2017": the dividend-tracker F column (cost basis pulled from
#    'Initial Buys') is flattened to plain cached values, with a few of them
#    nudged to their rounded cent amounts.
# ---------------------------------------------------------------------------
$fValues = @{
    2  = 2038.81
    3  = 1089.68
    4  = 946.8
    5  = 999.8
    6  = 1625.82
    7  = 1195.92
    8  = 1069.29
    9  = 1892.99
    10 = 948
    11 = 1062.5999999999999
    12 = 851.92
    13 = 1173.32
    14 = 2044.73
    15 = 999.4
    16 = 1356.18
    17 = 1078.6500000000001
    18 = 1042.68
}
foreach ($r in $fValues.Keys) {
    $ws2017.Range("F$r").Value = $fValues[$r]
}

# Row 17: January dividend for this holding was reversed out (S17) and the
# share count (D17) was corrected down slightly.
$ws2017.Range("D17").Value = 15.523999999999999
$ws2017.Range("S17").Value = 0

# ---------------------------------------------------------------------------
# 4. View-state: "2018" becomes the active/selected tab (activeCell D18),
#    "2017" keeps plain selection F3.
# ---------------------------------------------------------------------------
$ws2017.Range("F3").Select() | Out-Null
$ws2018.Activate() | Out-Null
$ws2018.Range("D18").Select() | Out-Null
